$wb = $excel.ActiveWorkbook

# The file "e6bb96d2-1476-4145-bc5c-0d402326e178.md" has been handed off again,
# moving its status from "In Translation" to "Ready for handoff" and bumping the
# "Latest Handoff"/"Latest Handoff Datetime" timestamps for each locale, as part
# of generating the handoff report.

# --- zh-cn sheet: row 3 corresponds to e6bb96d2-... .md ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("E3").Value = "2016-03-22 04:19:23"

# --- de-de sheet: row 3 corresponds to e6bb96d2-... .md ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("E3").Value = "2016-03-22 04:19:26"

# --- Overview sheet: row 3 rolls up the latest status/handoff date across locales ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-03-22 04:19:26"
